$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from the last existing
# header cell (AC1) onto the three new header cells before setting values.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-50) gets the team's season record appended.
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 80  # AD
    $ws.Cells.Item($r, 31).Value = 82  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
